# Auto-generated Excel COM-interop script
# Applies market-price / profit value updates to the Ultros profit tracker sheets
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 12319.685
$ws.Range("I64").Value = 5354
$ws.Range("J64").Value = 14177.2
$ws.Range("K64").Value = 5354
$ws.Range("L64").Value = 14177.2
$ws.Range("M64").Value = -5106
$ws.Range("N64").Value = -14673.2
$ws.Range("H67").Value = 12319.685
$ws.Range("I67").Value = 5354
$ws.Range("J67").Value = 14177.2
$ws.Range("K67").Value = 5354
$ws.Range("L67").Value = 14177.2
$ws.Range("M67").Value = -4496
$ws.Range("N67").Value = -15893.2
$ws.Range("H100").Value = 7086.7617
$ws.Range("I100").Value = 6567.25
$ws.Range("J100").Value = 7406.4614
$ws.Range("K100").Value = 6567.25
$ws.Range("L100").Value = 7406.4614
$ws.Range("M100").Value = -6026.25
$ws.Range("N100").Value = -8488.4614
$ws.Range("H134").Value = 38846.152
$ws.Range("J134").Value = 38846.152
$ws.Range("L134").Value = 38846.152
$ws.Range("N134").Value = -48986.152
$ws.Range("H137").Value = 3640.1072
$ws.Range("I137").Value = 2654.353
$ws.Range("J137").Value = 5163.5454
$ws.Range("K137").Value = 7963.059
$ws.Range("L137").Value = 15490.6362
$ws.Range("M137").Value = -5413.059
$ws.Range("N137").Value = -20590.6362

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 85429580
$ws.Range("J5").Value = 333334660
$ws.Range("L5").Value = 333334660
$ws.Range("N5").Value = -333334884
$ws.Range("H45").Value = 2856.4666
$ws.Range("I45").Value = 2881.9285
$ws.Range("K45").Value = 2881.9285
$ws.Range("M45").Value = -2504.9285
$ws.Range("H61").Value = 6726.28
$ws.Range("I61").Value = 1469.9412
$ws.Range("J61").Value = 17896
$ws.Range("K61").Value = 1469.9412
$ws.Range("L61").Value = 17896
$ws.Range("M61").Value = -1257.9412
$ws.Range("N61").Value = -18320
$ws.Range("H63").Value = 17849.9
$ws.Range("I63").Value = 15625
$ws.Range("J63").Value = 19333.166
$ws.Range("K63").Value = 15625
$ws.Range("L63").Value = 19333.166
$ws.Range("M63").Value = -14939
$ws.Range("N63").Value = -20705.166
$ws.Range("H66").Value = 17849.9
$ws.Range("I66").Value = 15625
$ws.Range("J66").Value = 19333.166
$ws.Range("K66").Value = 78125
$ws.Range("L66").Value = 96665.83
$ws.Range("M66").Value = -74693
$ws.Range("N66").Value = -103529.83
$ws.Range("H110").Value = 3068.6667
$ws.Range("I110").Value = 2733.7896
$ws.Range("K110").Value = 2733.7896
$ws.Range("M110").Value = -688.7896000000001
$ws.Range("H136").Value = 6726.28
$ws.Range("I136").Value = 1469.9412
$ws.Range("J136").Value = 17896
$ws.Range("K136").Value = 4409.8236
$ws.Range("L136").Value = 53688
$ws.Range("M136").Value = -1859.8236
$ws.Range("N136").Value = -58788

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 85429580
$ws.Range("J4").Value = 333334660
$ws.Range("L4").Value = 333334660
$ws.Range("N4").Value = -333334890
$ws.Range("H50").Value = 70000
$ws.Range("J50").Value = 70000
$ws.Range("L50").Value = 70000
$ws.Range("N50").Value = -71148
$ws.Range("H80").Value = 477.125
$ws.Range("I80").Value = 600
$ws.Range("J80").Value = 459.57144
$ws.Range("K80").Value = 600
$ws.Range("L80").Value = 459.57144
$ws.Range("M80").Value = 398
$ws.Range("N80").Value = -2455.57144
$ws.Range("H83").Value = 477.125
$ws.Range("I83").Value = 600
$ws.Range("J83").Value = 459.57144
$ws.Range("K83").Value = 3000
$ws.Range("L83").Value = 2297.8572
$ws.Range("M83").Value = 1992
$ws.Range("N83").Value = -12281.8572
$ws.Range("H86").Value = 17282520
$ws.Range("J86").Value = 6789.3
$ws.Range("L86").Value = 6789.3
$ws.Range("N86").Value = -9035.299999999999
$ws.Range("H89").Value = 17282520
$ws.Range("J89").Value = 6789.3
$ws.Range("L89").Value = 33946.5
$ws.Range("N89").Value = -45178.5
$ws.Range("H94").Value = 2050.6667
$ws.Range("I94").Value = 621.41174
$ws.Range("J94").Value = 8125
$ws.Range("K94").Value = 621.41174
$ws.Range("L94").Value = 8125
$ws.Range("M94").Value = -170.41174
$ws.Range("N94").Value = -9027
$ws.Range("H134").Value = 2055.3333
$ws.Range("I134").Value = 1687.375
$ws.Range("J134").Value = 4999
$ws.Range("K134").Value = 5062.125
$ws.Range("L134").Value = 14997
$ws.Range("M134").Value = -2527.125
$ws.Range("N134").Value = -20067

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1000
$ws.Range("I105").Value = 1000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1000
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 747
$ws.Range("N105").ClearContents()
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H141").Value = 49545.453
$ws.Range("J141").Value = 49545.453
$ws.Range("L141").Value = 49545.453
$ws.Range("N141").Value = -59905.453

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 540
$ws.Range("I8").Value = 540
$ws.Range("K8").Value = 1620
$ws.Range("M8").Value = -1481
$ws.Range("H12").Value = 517.43335
$ws.Range("J12").Value = 565.86957
$ws.Range("L12").Value = 1697.60871
$ws.Range("N12").Value = -2043.60871
$ws.Range("H50").Value = 1745.5714
$ws.Range("I50").Value = 294.9
$ws.Range("J50").Value = 5372.25
$ws.Range("K50").Value = 884.6999999999999
$ws.Range("L50").Value = 16116.75
$ws.Range("M50").Value = -403.6999999999999
$ws.Range("N50").Value = -17078.75
$ws.Range("H53").Value = 1745.5714
$ws.Range("I53").Value = 294.9
$ws.Range("J53").Value = 5372.25
$ws.Range("K53").Value = 884.6999999999999
$ws.Range("L53").Value = 16116.75
$ws.Range("M53").Value = -403.6999999999999
$ws.Range("N53").Value = -17078.75
$ws.Range("H93").Value = 5500
$ws.Range("J93").Value = 5000
$ws.Range("L93").Value = 15000
$ws.Range("N93").Value = -18744

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 4166826.2
$ws.Range("I2").Value = 56
$ws.Range("K2").Value = 56
$ws.Range("M2").Value = 57
$ws.Range("H97").Value = 17335
$ws.Range("I97").Value = 999.6667
$ws.Range("K97").Value = 999.6667
$ws.Range("M97").Value = -503.6667
$ws.Range("H102").Value = 2886.2144
$ws.Range("I102").Value = 2886.2144
$ws.Range("K102").Value = 2886.2144
$ws.Range("M102").Value = -1264.2144
$ws.Range("H132").Value = 6511.231
$ws.Range("I132").Value = 5821.095
$ws.Range("K132").Value = 17463.285
$ws.Range("M132").Value = -14933.285

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 13889425
$ws.Range("I16").Value = 17857304
$ws.Range("J16").Value = 1849.75
$ws.Range("K16").Value = 17857304
$ws.Range("L16").Value = 1849.75
$ws.Range("M16").Value = -17857134
$ws.Range("N16").Value = -2189.75
$ws.Range("H22").Value = 9091684
$ws.Range("I22").Value = 12987638
$ws.Range("K22").Value = 12987638
$ws.Range("M22").Value = -12987343
$ws.Range("H27").Value = 9091684
$ws.Range("I27").Value = 12987638
$ws.Range("K27").Value = 12987638
$ws.Range("M27").Value = -12987531
$ws.Range("H68").Value = 4187.125
$ws.Range("I68").Value = 4124.25
$ws.Range("K68").Value = 4124.25
$ws.Range("M68").Value = -3375.25
$ws.Range("H71").Value = 4187.125
$ws.Range("I71").Value = 4124.25
$ws.Range("K71").Value = 20621.25
$ws.Range("M71").Value = -16877.25
$ws.Range("H96").Value = 53998.668
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3006.3333
$ws.Range("I62").Value = 2220.4
$ws.Range("J62").Value = 3399.3
$ws.Range("K62").Value = 2220.4
$ws.Range("L62").Value = 3399.3
$ws.Range("M62").Value = -1596.4
$ws.Range("N62").Value = -4647.3
$ws.Range("H65").Value = 3006.3333
$ws.Range("I65").Value = 2220.4
$ws.Range("J65").Value = 3399.3
$ws.Range("K65").Value = 11102
$ws.Range("L65").Value = 16996.5
$ws.Range("M65").Value = -7982
$ws.Range("N65").Value = -23236.5
$ws.Range("H81").Value = 18521668
$ws.Range("I81").Value = 2300
$ws.Range("J81").Value = 37041036
$ws.Range("K81").Value = 4600
$ws.Range("L81").Value = 74082072
$ws.Range("M81").Value = -3539
$ws.Range("N81").Value = -74084194
$ws.Range("H84").Value = 18521668
$ws.Range("I84").Value = 2300
$ws.Range("J84").Value = 37041036
$ws.Range("K84").Value = 23000
$ws.Range("L84").Value = 370410360
$ws.Range("M84").Value = -17696
$ws.Range("N84").Value = -370420968
$ws.Range("H100").Value = 550.2727
$ws.Range("I100").Value = 571.8570999999999
$ws.Range("J100").Value = 512.5
$ws.Range("K100").Value = 1143.7142
$ws.Range("L100").Value = 1025
$ws.Range("M100").Value = -602.7141999999999
$ws.Range("N100").Value = -2107
$ws.Range("H126").Value = 1000.6667
$ws.Range("I126").Value = 952
$ws.Range("K126").Value = 2856
$ws.Range("M126").Value = -386

